$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect to make edits, then restore protection at the end.
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure footer (cell A38).
$disclosure = $ws.Range("A38").Value2
$disclosure = $disclosure.Replace("2021-04-30", "2021-05-03")
$ws.Range("A38").Value2 = $disclosure
$ws.Rows(38).AutoFit()

# Refresh the Weight (D) and Percent Change (E) columns for each holding row.
$ws.Range("D2").Value2 = 0.03624559449232355
$ws.Range("E2").Value2 = 0.0007745933384972759
$ws.Range("D3").Value2 = 0.02033816596015705
$ws.Range("E3").Value2 = 0.005090054815975131
$ws.Range("D4").Value2 = 0.01911059880805312
$ws.Range("E4").Value2 = 0.005075111652456332
$ws.Range("D5").Value2 = 0.03762640335228742
$ws.Range("E5").Value2 = 0.007407407407407307
$ws.Range("D6").Value2 = 0.0347407425443851
$ws.Range("E6").Value2 = -0.01106719367588926
$ws.Range("D7").Value2 = 0.01989666445073644
$ws.Range("E7").Value2 = -0.003463536655762933
$ws.Range("D8").Value2 = 0.03690256343974444
$ws.Range("E8").Value2 = 0.007718696397941649
$ws.Range("D9").Value2 = 0.02042392414443577
$ws.Range("E9").Value2 = -0.004498830304120904
$ws.Range("D10").Value2 = 0.02569937968757171
$ws.Range("E10").Value2 = 0.0001986294567484226
$ws.Range("D11").Value2 = 0.02355797740751682
$ws.Range("E11").Value2 = 0.008125677139761756
$ws.Range("D12").Value2 = 0.05669841097741513
$ws.Range("E12").Value2 = 0.007932692307692113
$ws.Range("D13").Value2 = 0.02490050136378484
$ws.Range("E13").Value2 = 0.007749077490774914
$ws.Range("D14").Value2 = 0.02769754538126583
$ws.Range("E14").Value2 = -0.001813510654375294
$ws.Range("D15").Value2 = 0.03303976979644734
$ws.Range("E15").Value2 = 0.002249524139124492
$ws.Range("D16").Value2 = 0.01913760242679327
$ws.Range("E16").Value2 = 0.03053040103492877
$ws.Range("D17").Value2 = 0.03075043464851166
$ws.Range("E17").Value2 = -0.004067065073041221
$ws.Range("D18").Value2 = 0.04224815692645118
$ws.Range("E18").Value2 = -0.0009182736455461971
$ws.Range("D19").Value2 = 0.1263143526479652
$ws.Range("E19").Value2 = 0.001327140013271499
$ws.Range("D20").Value2 = 0.008746824331046602
$ws.Range("E20").Value2 = 0.01986577181208071
$ws.Range("D21").Value2 = 0.01563157303940328
$ws.Range("E21").Value2 = 0.0001371553970648254
$ws.Range("D22").Value2 = 0.0166615900884366
$ws.Range("E22").Value2 = -0.001822916666666563
$ws.Range("D23").Value2 = 0.01634168143315402
$ws.Range("E23").Value2 = 0.004085801838610736
$ws.Range("D24").Value2 = 0.0217722163595273
$ws.Range("E24").Value2 = -0.001613065833249383
$ws.Range("D25").Value2 = 0.01208327711833307
$ws.Range("E25").Value2 = 0.01245400509482053
$ws.Range("D26").Value2 = 0.04169226012478735
$ws.Range("E26").Value2 = 0.009476584022038548
$ws.Range("D27").Value2 = 0.02394633947060758
$ws.Range("E27").Value2 = -0.0004902922141595267
$ws.Range("D28").Value2 = 0.0457532164296142
$ws.Range("E28").Value2 = 0.0004752851711029571
$ws.Range("D29").Value2 = 0.0554383782269771
$ws.Range("E29").Value2 = 0.003428983937917573
$ws.Range("D30").Value2 = 0.0134393284072403
$ws.Range("E30").Value2 = -0.005044136191677095
$ws.Range("D31").Value2 = 0.02065087705354481
$ws.Range("E31").Value2 = -0.001149425287356398
$ws.Range("D32").Value2 = 0.01406970210822714
$ws.Range("E32").Value2 = -0.01088435374149654
$ws.Range("D33").Value2 = 0.04193299559922689
$ws.Range("E33").Value2 = 0
$ws.Range("D34").Value2 = 0.01651095175402797
$ws.Range("E34").Value2 = 0.01287683684290242
$ws.Range("E35").Value2 = 0.002589080420605239

# Restore sheet protection.
$ws.Protect()
